# Consolidate split "word" + "space" text runs into a single run,
# e.g. <a:t>The</a:t><a:t> </a:t> -> <a:t>The </a:t>, leaving the
# trailing word run untouched.

$p = $ppt.ActivePresentation

# Slide 2: TextBox "The" + " " + "Moon" -> "The " + "Moon"
$s2 = $p.Slides.Item(2)
$tb2 = $s2.Shapes.Item(2).TextFrame.TextRange
$tb2.Characters(1, 4).Text = "The "

# Slide 3: Title "One" + " " + "More" -> "One " + "More"
$s3 = $p.Slides.Item(3)
$title3 = $s3.Shapes.Item(1).TextFrame.TextRange
$title3.Characters(1, 4).Text = "One "

# Slide 3: TextBox "The" + " " + "Moon" -> "The " + "Moon"
$tb3 = $s3.Shapes.Item(3).TextFrame.TextRange
$tb3.Characters(1, 4).Text = "The "
